{"js": "// Remove the trailing page-footer block from the end of the document:\n//   - a blank spacer paragraph\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ...\"\nconst body = context.document.body;\n\nconst results = body.search(\"Ver no Jupiter\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].paragraphs.load(\"items\");\n  await context.sync();\n\n  const hitPara = results.items[0].paragraphs.items[0];\n  const prevPara = hitPara.previous();   // blank spacer paragraph before it\n  const nextPara = hitPara.next();       // the \"\u00a9 2020 ...\" copyright paragraph\n  prevPara.load(\"text\");\n  nextPara.load(\"text\");\n  await context.sync();\n\n  // Delete from the end backwards so earlier deletions don't invalidate\n  // references to the paragraphs still to be removed.\n  nextPara.delete();\n  hitPara.delete();\n  prevPara.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The page footer block (\"Ver no Jupiter Salvar em pdf Salvar em docx\" plus the\n# \"(c) 2020 ... Jekyll ...\" copyright line, and the blank spacer paragraph right\n# before them) is being removed from the end of the document.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Ver no Jupiter*\") {\n        $target = $p\n    }\n}\n\nif ($target -ne $null) {\n    $startPara = $target.Previous()   # blank spacer paragraph\n    $endPara = $target.Next()         # the \"(c) 2020 ...\" copyright paragraph\n    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rng.Delete()\n}\n"}
